$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the question text in A4: remove the trailing clause about the keyboard
# appearing below, per the authored edit.
$ws.Range("A4").Value = "Czy zauważyłeś coś charakterstycznego w słowach i obrazach, które pojawiały się z określonymi postaciami z kreskówek?`n`nWpisz swoją odpowiedź za pomocą klawiatury.`n`nNaciśnij Enter by przejść do następnego pytania.`n`nTwoja odpowiedź musi mieć minimum 20 znaków."

# Remove sheet protection (sheetProtection element no longer present in target).
$ws.Unprotect()
